# "Motor speed Finally the same RPM"
#
# 1. Rename the worksheet to reflect that it now tracks motor speed vs ticks.
# 2. Re-enter the "Average" formulas as single range assignments so Excel
#    stores each block (rows 3-10 and rows 14-21) as one shared formula
#    group, same as the resaved workbook.
# 3. Leave the sheet's saved selection on the merged header cell A1:A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Motor Speed vs Ticks"

$ws.Range("L3:L10").Formula = "=AVERAGE(B3:K3)"
$ws.Range("L14:L21").Formula = "=AVERAGE(B14:K14)"

$ws.Range("A1:A2").Select() | Out-Null
